# Fruta / hortaliza, semanal
# Insert 2 new weekly-price rows before the current row 100 (pushing existing
# rows 100-109 down to 102-111) and populate the two new rows with the
# latest week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows right above row 100; this shifts old rows 100-109
# down to 102-111, exactly matching the diff.
$ws.Rows("100:101").Insert()

# --- New row 100 ---
$ws.Range("A100").Value = 4
$ws.Range("B100").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C100").Value = "Los Lagos"
$ws.Range("D100").Value = 44491
$ws.Range("E100").Value = 10
$ws.Range("F100").Value = "Fruta"
$ws.Range("G100").Value = 100101
$ws.Range("H100").Value = "Berries"
$ws.Range("I100").Value = 100112025
$ws.Range("J100").Value = "Frutilla"
$ws.Range("K100").Value = "Sin especificar"
$ws.Range("L100").Value = "Especial"
$ws.Range("M100").Value = 300
$ws.Range("N100").Value = 12500
$ws.Range("O100").Value = 12500
$ws.Range("P100").Value = 12500
$ws.Range("Q100").Value = "$/bandeja 7 kilos"
$ws.Range("R100").Value = "Provincia de Melipilla"
$ws.Range("S100").Value = 1786
$ws.Range("T100").Value = 7

# --- New row 101 ---
$ws.Range("A101").Value = 4
$ws.Range("B101").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C101").Value = "Los Lagos"
$ws.Range("D101").Value = 44491
$ws.Range("E101").Value = 10
$ws.Range("F101").Value = "Fruta"
$ws.Range("G101").Value = 100101
$ws.Range("H101").Value = "Berries"
$ws.Range("I101").Value = 100112025
$ws.Range("J101").Value = "Frutilla"
$ws.Range("K101").Value = "Sin especificar"
$ws.Range("L101").Value = "Primera"
$ws.Range("M101").Value = 600
$ws.Range("N101").Value = 10000
$ws.Range("O101").Value = 10500
$ws.Range("P101").Value = 10250
$ws.Range("Q101").Value = "$/bandeja 7 kilos"
$ws.Range("R101").Value = "Provincia de Melipilla"
$ws.Range("S101").Value = 1464
$ws.Range("T101").Value = 7

# Make sure the sheet dimension / used range reflects the new last row.
$ws.Range("A1:T111").Select()
